$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record arrived. Insert a new data row right above the
# current row 24 (the first "Albahaca" data row after the older, out-of-week
# batch that occupies rows 2-23), pushing the existing rows 24-70 down to
# 25-71, and populate the freshly inserted row 24 with the new record.
$ws.Rows.Item(24).Insert()

$ws.Cells.Item(24, 1).Value = 8
$ws.Cells.Item(24, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(24, 3).Value = "Coquimbo"
$ws.Cells.Item(24, 4).Value = 44544
$ws.Cells.Item(24, 5).Value = 4
$ws.Cells.Item(24, 6).Value = 100112052
$ws.Cells.Item(24, 7).Value = "Albahaca"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 800
$ws.Cells.Item(24, 11).Value = 3000
$ws.Cells.Item(24, 12).Value = 4000
$ws.Cells.Item(24, 13).Value = 3500
$ws.Cells.Item(24, 14).Value = "`$/paquete"
$ws.Cells.Item(24, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(24, 16).Value = 3500
$ws.Cells.Item(24, 17).Value = 1
$ws.Cells.Item(24, 18).Value = "Hortaliza"
